$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 35 / 36 text values (historical reorder of entries)
$ws.Range("C35").Value = "ajax funktio luotu, axios implementation myöhemmin"
$ws.Range("C36").Value = "perustoimintojen alustava viimeistely, date input=> fetch=> datan manipulointi=> datan esittely=> perus css"

# Push the totals row (old row 37: date label + SUM formula) down by 5 rows so
# 5 new log entries can be inserted above it, at rows 37-41.
$ws.Rows("37").Insert()
$ws.Rows("37").Insert()
$ws.Rows("37").Insert()
$ws.Rows("37").Insert()
$ws.Rows("37").Insert()

# New row 37
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A37").PasteSpecial(-4122) | Out-Null
$ws.Range("A37").Value = 44534
$ws.Range("B37").Value = 1
$ws.Range("C37").Value = "pientä css tuunausta ja testifunktioiden siivousta"

# New row 38
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A38").PasteSpecial(-4122) | Out-Null
$ws.Range("A38").Value = 44535
$ws.Range("B38").Value = 2
$ws.Range("C38").Value = "rakenteen refaktorin, Loading komponentti ja turhat statet pois, 'ylimääräiset' useEffect hookit pois"

# New row 39
$ws.Range("B39").Value = 2
$ws.Range("C39").Value = "Docker image(production testi) ja testiajoa kontissa"

# New row 40
$ws.Range("B40").Value = 1
$ws.Range("C40").Value = "CI/CD push to github main => uusi kuva dockerhubiin"

# New row 41
$ws.Range("B41").Value = 1
$ws.Range("C41").Value = "Heroku app, workflow push github main => uusi appi Herokuun"

# Row 42 (the old totals row, now shifted): formula range needs to extend to B41
$ws.Range("B42").Formula = "=SUM(B2:B41)"

$ws.Range("C41").Select()
$excel.ActiveWindow.ScrollRow = 25

$wb.Save()
